$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (H1) onto the two new
# header cells (I1, J1) so they match the rest of the header row, then set
# their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-10 for new columns I (9) and J (10)
$data = @{
    2  = @(1, 1)
    3  = @(1, 3)
    4  = @(3, 7)
    5  = @(6, 7)
    6  = @(3, 6)
    7  = @(13, 13)
    8  = @(4, 8)
    9  = @(1, 2)
    10 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
